$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells (new shared strings: "diferencia2", "merma2") ---
$ws.Range("G1").Value = "diferencia2"
$ws.Range("H1").Value = "merma2"

# --- Column G: diferencia2 = (B/D) - 1 ---
# G2 stays a standalone formula (matches columns E/F, whose row-2 formula is
# also not part of the shared-formula group); G3:G22 become one shared formula.
$ws.Range("G2").Formula = "=SUM((B2/D2),-1)"
$ws.Range("G3:G22").Formula = "=SUM((B3/D3),-1)"

# --- Column H: merma2 = literal copy of column E's computed values (no formula) ---
for ($r = 2; $r -le 22; $r++) {
    $eVal = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($r, 8).Value = $eVal
}

# --- Number format "0.0" for the new columns (creates numFmtId 164 / cellXfs index 2) ---
$ws.Range("G1:H22").NumberFormat = "0.0"

# --- Column widths (closest achievable widths via the ColumnWidth->stored-width model) ---
$ws.Columns.Item(7).ColumnWidth = 17.166666666666668
$ws.Columns.Item(8).ColumnWidth = 13.5

# --- Selection moved to G2 (matches the saved view state) ---
$ws.Range("G2").Select() | Out-Null
